$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'45.912.79"
$ws.Range("E2").Value = "  -0.27%  "

$ws.Range("D3").Value = "'2.603.17"
$ws.Range("E3").Value = "  +0.82%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").Value = "'308.08"
$ws.Range("E5").Value = "  +0.78%  "

$ws.Range("D6").Value = "'98.85"
$ws.Range("E6").Value = "  -0.76%  "

$ws.Range("D7").Value = "'0.594"
$ws.Range("E7").Value = "  -0.20%  "

$ws.Range("E8").Value = "  +0.04%  "

$ws.Range("D9").Value = "'0.579"
$ws.Range("E9").Value = "  +1.05%  "

$ws.Range("D10").Value = "'38.72"
$ws.Range("E10").Value = "  +0.73%  "

$ws.Range("B11").Value = "OKB"
$ws.Range("C11").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D11").Value = "'54.18"
$ws.Range("E11").Value = "  -0.62%  "

$ws.Range("B12").Value = "Dogecoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D12").Value = "'0.0841"
$ws.Range("E12").Value = "  +0.32%  "

$ws.Range("D13").Value = "'8.08"
$ws.Range("E13").Value = "  -2.68%  "

$ws.Range("D14").Value = "'3.002.26"
$ws.Range("E14").Value = "  +0.80%  "

$ws.Range("D16").Value = "'2.606.72"
$ws.Range("E16").Value = "  +0.20%  "

$ws.Range("D17").Value = "'0.910"
$ws.Range("E17").Value = "  +1.12%  "

$ws.Range("D18").Value = "'14.78"
$ws.Range("E18").Value = "  -0.15%  "

$ws.Range("D19").Value = "'46.125.81"
$ws.Range("E19").Value = "  -0.09%  "

$ws.Range("D20").Value = "'0.0000101"
$ws.Range("E20").Value = "  +0.35%  "

$ws.Range("D21").Value = "'6.75"
$ws.Range("E21").Value = "  +1.96%  "

$ws.Range("D22").Value = "'12.64"
$ws.Range("E22").Value = "  -2.15%  "

$ws.Range("D23").Value = "'290.77"
$ws.Range("E23").Value = "  +14.91%  "

$ws.Range("D24").Value = "'72.52"
$ws.Range("E24").Value = "  +2.11%  "

$ws.Range("D25").Value = "'3.02"
$ws.Range("E25").Value = "  +1.23%  "

$ws.Range("E26").Value = "  +2.64%  "

$ws.Range("D27").Value = "'29.57"
$ws.Range("E27").Value = "  +4.78%  "

$ws.Range("E28").Value = "  +0.08%  "

$ws.Range("D29").Value = "'4.06"
$ws.Range("E29").Value = "  +1.20%  "

$ws.Range("D30").Value = "'10.73"
$ws.Range("E30").Value = "  +3.04%  "

$ws.Range("D31").Value = "'38.63"
$ws.Range("E31").Value = "  -2.36%  "

$ws.Range("E32").Value = "  -2.51%  "

$ws.Range("D33").Value = "'6.24"
$ws.Range("E33").Value = "  +3.38%  "

$ws.Range("B34").Value = "LidoDAOToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D34").Value = "'3.63"
$ws.Range("E34").Value = "  -1.04%  "

$ws.Range("B35").Value = "Monero"
$ws.Range("C35").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D35").Value = "'159.66"
$ws.Range("E35").Value = "  +4.63%  "

$ws.Range("D36").Value = "'2.23"
$ws.Range("E36").Value = "  -2.13%  "

$ws.Range("D37").Value = "'0.0839"
$ws.Range("E37").Value = "  +2.07%  "

$ws.Range("E38").Value = "  -3.77%  "

$ws.Range("E39").Value = "  +4.39%  "

$ws.Range("E40").Value = "  +1.16%  "

$ws.Range("D41").Value = "'15.67"
$ws.Range("E41").Value = "  -2.70%  "

$ws.Range("D42").Value = "'0.0329"
$ws.Range("E42").Value = "  +2.98%  "

$ws.Range("D43").Value = "'3.53"
$ws.Range("E43").Value = "  -1.85%  "

$ws.Range("D44").Value = "'21.42"
$ws.Range("E44").Value = "  +8.03%  "

$ws.Range("D45").Value = "'4.00"
$ws.Range("E45").Value = "  -4.05%  "

$ws.Range("D46").Value = "'2.110.55"
$ws.Range("E46").Value = "  +2.86%  "

$ws.Range("D47").Value = "'95.85"
$ws.Range("E47").Value = "  +5.61%  "

$ws.Range("D48").Value = "'0.999"
$ws.Range("E48").Value = "  +0.02%  "

$ws.Range("D49").Value = "'9.36"
$ws.Range("E49").Value = "  +1.07%  "

$ws.Range("D50").Value = "'109.31"
$ws.Range("E50").Value = "  +1.03%  "

$ws.Range("D51").Value = "'2.868.64"
$ws.Range("E51").Value = "  +1.03%  "
